# Update "想去人数" (interested-count) figures that changed between
# the previous and newly generated data snapshot.
#
# Sheet "展览" (Exhibition):
#   F5:  11085 -> 11086
#   F10: 10975 -> 10977
#   F15: 5497  -> 5499
#   F17: 3421  -> 3422
#
# Sheet "全部类型" (All types) mirrors the same events:
#   F7:  11085 -> 11086
#   F12: 10975 -> 10977
#   F17: 5497  -> 5499
#   F19: 3421  -> 3422

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 11086
$wsExhibit.Range("F10").Value = 10977
$wsExhibit.Range("F15").Value = 5499
$wsExhibit.Range("F17").Value = 3422

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 11086
$wsAll.Range("F12").Value = 10977
$wsAll.Range("F17").Value = 5499
$wsAll.Range("F19").Value = 3422
